# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.859.91'
$ws.Range("E2").Value = '  -0.24%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.852.09'
$ws.Range("E3").Value = '  +1.62%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '694.62'
$ws.Range("E5").Value = '  -1.13%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.55'
$ws.Range("E6").Value = '  -0.21%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.850.59'
$ws.Range("E7").Value = '  +1.58%  '

# Row 8
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("E9").Value = '  +0.05%  '

# Row 10
$ws.Range("E10").Value = '  -0.40%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.23'
$ws.Range("E11").Value = '  -3.14%  '

# Row 12
$ws.Range("E12").Value = '  -0.42%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000256'
$ws.Range("E13").Value = '  +0.41%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.21'
$ws.Range("E14").Value = '  +0.28%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.504.00'
$ws.Range("E15").Value = '  +1.69%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.858.96'
$ws.Range("E16").Value = '  +1.62%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '70.953.64'
$ws.Range("E17").Value = '  -0.12%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.18'
$ws.Range("E18").Value = '  -0.06%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.40'
$ws.Range("E19").Value = '  -2.67%  '

# Row 20
$ws.Range("E20").Value = '  -0.12%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '497.41'
$ws.Range("E21").Value = '  +3.09%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.68'
$ws.Range("E22").Value = '  -4.71%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.719'
$ws.Range("E23").Value = '  +0.56%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.87'
$ws.Range("E24").Value = '  +1.17%  '

# Row 25
$ws.Range("E25").Value = '  +1.86%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.66'
$ws.Range("E26").Value = '  +1.22%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.20'
$ws.Range("E27").Value = '  -1.31%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.12'
$ws.Range("E28").Value = '  -2.89%  '

# Row 29
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.03%  '

# Row 30
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.14'
$ws.Range("E30").Value = '  +1.06%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.53'
$ws.Range("E31").Value = '  -0.29%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.26'
$ws.Range("E32").Value = '  -2.07%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.50'
$ws.Range("E33").Value = '  -0.04%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.181'
$ws.Range("E34").Value = '  +1.93%  '

# Row 35
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.20'
$ws.Range("E35").Value = '  +0.21%  '

# Row 36
$ws.Range("B36").Value = 'RenzoRestakedETH'
$ws.Range("C36").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.809.34'
$ws.Range("E36").Value = '  +1.85%  '

# Row 37
$ws.Range("E37").Value = '  -0.05%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.103'
$ws.Range("E38").Value = '  +0.91%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.38'
$ws.Range("E39").Value = '  +6.87%  '

# Row 40
$ws.Range("E40").Value = '  +9.08%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.38'
$ws.Range("E41").Value = '  -2.04%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.02'
$ws.Range("E42").Value = '  +0.73%  '

# Row 43
$ws.Range("E43").Value = '  +0.01%  '

# Row 44
$ws.Range("E44").Value = '  +0.20%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000312'
$ws.Range("E45").Value = '  -3.85%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '163.95'
$ws.Range("E46").Value = '  +1.68%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '49.15'
$ws.Range("E47").Value = '  +0.25%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.300'
$ws.Range("E48").Value = '  +0.14%  '

# Row 49
$ws.Range("B49").Value = 'Arweave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '43.52'
$ws.Range("E49").Value = '  -5.67%  '

# Row 50
$ws.Range("B50").Value = 'ONDO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.38'
$ws.Range("E50").Value = '  -3.22%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.62'
$ws.Range("E51").Value = '  +1.43%  '
